$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: XauUsd
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("XauUsd")

$ws1.Range("A2").Value = 245
$ws1.Range("B2").Value = 245
$ws1.Range("D2").Value = 245
$ws1.Range("E2").Value = 245

$ws1.Range("A3").Value = 365.05
$ws1.Range("D3").Value = 365.05
$ws1.Range("B3").Value = 185.5
$ws1.Range("E3").Value = 185.5

$ws1.Range("A5").Value = 61.51
$ws1.Range("D5").Value = 61.51
$ws1.Range("B5").Value = 64.27
$ws1.Range("E5").Value = 64.27

$ws1.Range("A6").Value = 54.41
$ws1.Range("D6").Value = 54.41
$ws1.Range("B6").Value = 51.29
$ws1.Range("E6").Value = 51.29

$ws1.Range("A12").Formula = "=B3-E3"

# ---------------------------------------------------------------------------
# Sheet 2: UsdJpy
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("UsdJpy")

$ws2.Range("A2").Value = 193
$ws2.Range("D2").Value = 193
$ws2.Range("B2").Value = 193
$ws2.Range("E2").Value = 193

$ws2.Range("A3").Value = 825.4
$ws2.Range("D3").Value = 825.4
$ws2.Range("B3").Value = 283.25
$ws2.Range("E3").Value = 283.25

$ws2.Range("A5").Value = 51.84
$ws2.Range("D5").Value = 51.84
$ws2.Range("B5").Value = 52.23
$ws2.Range("E5").Value = 52.23

$ws2.Range("A6").Value = 48.32
$ws2.Range("D6").Value = 48.32
$ws2.Range("B6").Value = 57.31
$ws2.Range("E6").Value = 57.31

$ws2.Range("A12").Formula = "=B3-E3"

# ---------------------------------------------------------------------------
# Sheet 3: UsdChf
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("UsdChf")

$ws3.Range("A2").Value = 105
$ws3.Range("D2").Value = 105
$ws3.Range("B2").Value = 105
$ws3.Range("E2").Value = 105

$ws3.Range("A3").Value = 378.4
$ws3.Range("D3").Value = 378.4
$ws3.Range("B3").Value = 142.30000000000001
$ws3.Range("E3").Value = 142.30000000000001

$ws3.Range("A5").Value = 41.93
$ws3.Range("D5").Value = 41.93
$ws3.Range("B5").Value = 45.5
$ws3.Range("E5").Value = 45.5

$ws3.Range("A6").Value = 42.4
$ws3.Range("D6").Value = 42.4
$ws3.Range("B6").Value = 47.83
$ws3.Range("E6").Value = 47.83

$ws3.Range("A12").Formula = "=B3-E3"

# ---------------------------------------------------------------------------
# Sheet 4: UsdCad
# ---------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("UsdCad")

$ws4.Range("A2").Value = 250
$ws4.Range("D2").Value = 250
$ws4.Range("B2").Value = 250
$ws4.Range("E2").Value = 250

$ws4.Range("A3").Value = 615.79999999999995
$ws4.Range("D3").Value = 615.79999999999995
$ws4.Range("B3").Value = 269.05
$ws4.Range("E3").Value = 269.05

$ws4.Range("A5").Value = 51.64
$ws4.Range("D5").Value = 51.64
$ws4.Range("B5").Value = 51.71
$ws4.Range("E5").Value = 51.71

$ws4.Range("A6").Value = 48.18
$ws4.Range("D6").Value = 48.18
$ws4.Range("B6").Value = 60.23
$ws4.Range("E6").Value = 60.23

$ws4.Range("A12").Formula = "=B3-E3"

# ---------------------------------------------------------------------------
# Selections per-sheet: all sheets now select D1:E8 with active cell D1
# ---------------------------------------------------------------------------
$ws1.Range("D1:E8").Select()
$ws2.Range("D1:E8").Select()
$ws3.Range("D1:E8").Select()

# Active sheet moves from UsdChf (index 3) to UsdCad (index 4)
$ws4.Activate()
$ws4.Range("D1:E8").Select()
